# Edit script: Update EC database table with full period history (ascending order)
# and append a new row for period 2508 (part 1 of new account statement), per commit:
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 125 for the new period (2508) ---------------------
$ws.Rows.Item(125).Insert()

# The inserted row copies the format of the row above (old "last row" style that
# used to belong to row 124). Move that special bottom-border formatting down to
# the new last row (125), and restore row 124 to the regular interior-row style
# taken from row 123.
$ws.Range("B124:J124").Copy()
$ws.Range("B125:J125").PasteSpecial(-4122)

$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Rewrite the "Periodo Mora" column (E16:E125) in ascending order -------
# Previously the table listed periods descending from 2507 (row16) down to 1607
# (row124). The refreshed database now lists them ascending from 1607 (row16) up
# to 2507 (row124), with the brand-new period 2508 appended as row125.
$periods = @(
  "1607","1608","1609","1610","1611","1612",
  "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
  "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
  "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
  "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
  "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
  "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
  "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
  "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
  "2501","2502","2503","2504","2505","2506","2507",
  "2508"
)

$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 5).Value2 = $p
    $row = $row + 1
}

# Fill in the rest of the new row 125 (B,C,D,F,G columns) matching the pattern
# used by every other data row in the table.
$ws.Cells.Item(125, 2).Value2 = "CC"
$ws.Cells.Item(125, 3).Value2 = "80092270"
$ws.Cells.Item(125, 4).Value2 = "JAVIER ANDRES QUEVEDO CORREA"
$ws.Cells.Item(125, 6).Value2 = 100000
$ws.Cells.Item(125, 7).Value2 = 2500000

# --- 3. Update the summary figures at the top of the sheet --------------------
# Total amount in arrears (VALOR MORA) grew from 10,900,000 to 11,000,000 ...
$ws.Cells.Item(11, 5).Value2 = 11000000

# ... and the period count grew from 109 to 110.
$ws.Cells.Item(13, 6).Value2 = 110
